$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Rows 3-8: columns D, E, F change to 2, 2, 2 (supply data for first loop block)
foreach ($r in 3..8) {
    $ws.Cells.Item($r, 4).Value = 2   # D
    $ws.Cells.Item($r, 5).Value = 2   # E
    $ws.Cells.Item($r, 6).Value = 2   # F
}

# Rows 9-14: column F changes from 8 to 1 (demand data for second loop block)
foreach ($r in 9..14) {
    $ws.Cells.Item($r, 6).Value = 1   # F
}

# Update the active cell selection to G12
$ws.Range("G12").Select()
